$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 3.080253
$ws.Range("N2").Value = 9.240759000000001
$ws.Range("O2").Value = 0.9515444181830145
$ws.Range("P2").Value = 0.9515444181830145
$ws.Range("Q2").Value = 1.478278100013
$ws.Range("R2").Value = 13.304502900117
$ws.Range("S2").Value = 0.9515444181830145
$ws.Range("T2").Value = 0.9515444181830145

# Row 3
$ws.Range("O3").Value = 0.002860371193349785
$ws.Range("P3").Value = 0.002860371193349786
$ws.Range("Q3").Value = 0.004443748512666667
$ws.Range("R3").Value = 0.03999373661400001
$ws.Range("S3").Value = 0.002860371193349785
$ws.Range("T3").Value = 0.002860371193349786

# Row 4
$ws.Range("O4").Value = 0.04559521062363567
$ws.Range("P4").Value = 0.04559521062363567
$ws.Range("R4").Value = 0.6375126587700001
$ws.Range("S4").Value = 0.04559521062363567
$ws.Range("T4").Value = 0.04559521062363567
